$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Both A1 and A2 now hold the same text, replacing the old "Headers" /
# flight-related strings.
$ws.Range("A1").Value = "Rahul Shetty Academy"
$ws.Range("A2").Value = "Rahul Shetty Academy"

# Apply the same formatting to both cells: small Consolas font in dark grey,
# not bold (A1 was previously bold Calibri; normalize it away).
$rng = $ws.Range("A1:A2")
$rng.Font.Name = "Consolas"
$rng.Font.Family = 3
$rng.Font.Size = 7
$rng.Font.Bold = $false
$rng.Font.Color = 2236962
